$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.94
$ws.Range("G2").Value = 3.5
$ws.Range("H2").Value = 2.36
$ws.Range("I2").Value = 2.76
$ws.Range("N2").Value = 3.5
$ws.Range("O2").Value = 1.31
$ws.Range("Q2").Value = 1.94
$ws.Range("T2").Value = 1.73
$ws.Range("V2").Value = 1.58
$ws.Range("W2").Value = 1.4
$ws.Range("Y2").Value = 13
$ws.Range("Z2").Value = 20
$ws.Range("AA2").Value = 44
$ws.Range("AD2").Value = 14.5
$ws.Range("AE2").Value = 34
$ws.Range("AI2").Value = 50
$ws.Range("AO2").Value = 27
$ws.Range("G3").Value = 2.58
$ws.Range("H3").Value = 2.78
$ws.Range("L3").Value = 1.32
$ws.Range("W3").Value = 1.63
$ws.Range("H4").Value = 2.86
$ws.Range("T4").Value = 1.64
$ws.Range("F5").Value = 8.6
$ws.Range("G5").Value = 8.800000000000001
$ws.Range("H5").Value = 1.43
$ws.Range("I5").Value = 1.44
$ws.Range("K5").Value = 5.3
$ws.Range("L5").Value = 1.34
$ws.Range("N5").Value = 4.8
$ws.Range("Q5").Value = 1.73
$ws.Range("R5").Value = 1.51
$ws.Range("S5").Value = 2.86
$ws.Range("T5").Value = 1.96
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 3.25
$ws.Range("W5").Value = 1.12
$ws.Range("Z5").Value = 8.6
$ws.Range("AG5").Value = 32
$ws.Range("AJ5").Value = 290
$ws.Range("G6").Value = 5
$ws.Range("J6").Value = 4.2
$ws.Range("M6").Value = 1.03
$ws.Range("N6").Value = 5.4
$ws.Range("O6").Value = 1.19
$ws.Range("Q6").Value = 1.58
$ws.Range("R6").Value = 1.61
$ws.Range("S6").Value = 2.36
$ws.Range("T6").Value = 1.6
$ws.Range("U6").Value = 2.44
$ws.Range("W6").Value = 1.25
$ws.Range("X6").Value = 25
$ws.Range("Y6").Value = 13
$ws.Range("Z6").Value = 14
$ws.Range("AA6").Value = 20
$ws.Range("AB6").Value = 970
$ws.Range("AC6").Value = 10
$ws.Range("AD6").Value = 11
$ws.Range("AE6").Value = 970
$ws.Range("AF6").Value = 40
$ws.Range("AG6").Value = 970
$ws.Range("AH6").Value = 16
$ws.Range("AI6").Value = 27
$ws.Range("AJ6").Value = 110
$ws.Range("AK6").Value = 50
$ws.Range("AL6").Value = 48
$ws.Range("AM6").Value = 65
$ws.Range("AN6").Value = 40
$ws.Range("AO6").Value = 7.8
$ws.Range("F7").Value = 4.4
$ws.Range("G7").Value = 4.5
$ws.Range("H7").Value = 1.87
$ws.Range("I7").Value = 1.89
$ws.Range("L7").Value = 1.33
$ws.Range("Q7").Value = 1.75
$ws.Range("S7").Value = 2.86
$ws.Range("T7").Value = 1.71
$ws.Range("V7").Value = 2.12
$ws.Range("W7").Value = 1.28
$ws.Range("X7").Value = 19.5
$ws.Range("Z7").Value = 12.5
$ws.Range("AA7").Value = 20
$ws.Range("AC7").Value = 9.199999999999999
$ws.Range("AD7").Value = 10
$ws.Range("AE7").Value = 17.5
$ws.Range("AF7").Value = 36
$ws.Range("AG7").Value = 17
$ws.Range("AH7").Value = 17
$ws.Range("AI7").Value = 29
$ws.Range("AJ7").Value = 95
$ws.Range("AK7").Value = 48
$ws.Range("AL7").Value = 50
$ws.Range("AN7").Value = 42
$ws.Range("AO7").Value = 10
$ws.Range("F8").Value = 1.04
$ws.Range("G8").Value = 1000
$ws.Range("H8").Value = 1.04
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 1.02
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 1.01
$ws.Range("M8").Value = 1.01
$ws.Range("N8").Value = 1.01
$ws.Range("O8").Value = 1.28
$ws.Range("P8").Value = 1.08
$ws.Range("Q8").Value = 1.28
$ws.Range("R8").Value = 1.08
$ws.Range("S8").Value = 1.28
$ws.Range("T8").Value = 1.01
$ws.Range("U8").Value = 1.01
$ws.Range("V8").Value = 1.01
$ws.Range("W8").Value = 1.01
$ws.Range("X8").Value = 1000
$ws.Range("Y8").Value = 1000
$ws.Range("Z8").Value = 1000
$ws.Range("AA8").Value = 1000
$ws.Range("AB8").Value = 1000
$ws.Range("AC8").Value = 1000
$ws.Range("AD8").Value = 1000
$ws.Range("AE8").Value = 1000
$ws.Range("AF8").Value = 1000
$ws.Range("AG8").Value = 1000
$ws.Range("AH8").Value = 1000
$ws.Range("AI8").Value = 1000
$ws.Range("AJ8").Value = 1000
$ws.Range("AK8").Value = 1000
$ws.Range("AL8").Value = 1000
$ws.Range("AM8").Value = 1000
$ws.Range("AN8").Value = 1000
$ws.Range("AO8").Value = 1000
$ws.Range("F9").Value = 2.38
$ws.Range("K9").Value = 3.45
$ws.Range("L9").Value = 1.4
$ws.Range("P9").Value = 1.89
$ws.Range("Q9").Value = 2.06
$ws.Range("S9").Value = 3.75
$ws.Range("V9").Value = 1.4
$ws.Range("W9").Value = 1.71
$ws.Range("X9").Value = 13.5
$ws.Range("Y9").Value = 12.5
$ws.Range("AB9").Value = 9.800000000000001
$ws.Range("AD9").Value = 14
$ws.Range("AE9").Value = 40
$ws.Range("AH9").Value = 18
$ws.Range("AL9").Value = 40
$ws.Range("AM9").Value = 100
$ws.Range("F10").Value = 13
$ws.Range("G10").Value = 13.5
$ws.Range("L10").Value = 1.31
$ws.Range("S10").Value = 2.88
$ws.Range("V10").Value = 4.2
$ws.Range("W10").Value = 1.08
$ws.Range("X10").Value = 20
$ws.Range("Y10").Value = 8.4
$ws.Range("AC10").Value = 13
$ws.Range("AF10").Value = 130
$ws.Range("AG10").Value = 46
$ws.Range("AI10").Value = 42
$ws.Range("AJ10").Value = 1000
$ws.Range("AK10").Value = 260
$ws.Range("AM10").Value = 230
$ws.Range("AN10").Value = 370
$ws.Range("AO10").Value = 5.2
$ws.Range("G11").Value = 2.74
$ws.Range("H11").Value = 2.88
$ws.Range("I11").Value = 2.9
$ws.Range("L11").Value = 1.38
$ws.Range("P11").Value = 2
$ws.Range("S11").Value = 3.5
$ws.Range("V11").Value = 1.52
$ws.Range("W11").Value = 1.57
$ws.Range("X11").Value = 14
$ws.Range("Z11").Value = 19
$ws.Range("AA11").Value = 44
$ws.Range("AD11").Value = 12.5
$ws.Range("AE11").Value = 30
$ws.Range("AF11").Value = 17.5
$ws.Range("AG11").Value = 12
$ws.Range("AH11").Value = 16.5
$ws.Range("AI11").Value = 44
$ws.Range("AK11").Value = 28
$ws.Range("AM11").Value = 85
$ws.Range("AO11").Value = 26
$ws.Range("H12").Value = 1.99
$ws.Range("L12").Value = 1.3
$ws.Range("O12").Value = 1.24
$ws.Range("R12").Value = 1.54
$ws.Range("S12").Value = 2.78
$ws.Range("U12").Value = 2.42
$ws.Range("V12").Value = 2
$ws.Range("W12").Value = 1.32
$ws.Range("X12").Value = 18.5
$ws.Range("Y12").Value = 12
$ws.Range("Z12").Value = 13.5
$ws.Range("AA12").Value = 23
$ws.Range("AB12").Value = 18.5
$ws.Range("AC12").Value = 8.800000000000001
$ws.Range("AD12").Value = 10
$ws.Range("AE12").Value = 18
$ws.Range("AF12").Value = 32
$ws.Range("AG12").Value = 16
$ws.Range("AH12").Value = 16
$ws.Range("AI12").Value = 28
$ws.Range("AK12").Value = 42
$ws.Range("AL12").Value = 46
$ws.Range("AN12").Value = 34
$ws.Range("F13").Value = 4.8
$ws.Range("G13").Value = 4.9
$ws.Range("H13").Value = 1.79
$ws.Range("I13").Value = 1.81
$ws.Range("J13").Value = 4.2
$ws.Range("K13").Value = 4.3
$ws.Range("L13").Value = 1.25
$ws.Range("M13").Value = 1.04
$ws.Range("N13").Value = 5.9
$ws.Range("O13").Value = 1.19
$ws.Range("P13").Value = 2.6
$ws.Range("Q13").Value = 1.59
$ws.Range("R13").Value = 1.65
$ws.Range("S13").Value = 2.42
$ws.Range("V13").Value = 2.26
$ws.Range("W13").Value = 1.25
$ws.Range("X13").Value = 25
$ws.Range("AA13").Value = 19.5
$ws.Range("AB13").Value = 25
$ws.Range("AC13").Value = 9.800000000000001
$ws.Range("AD13").Value = 9.800000000000001
$ws.Range("AF13").Value = 40
$ws.Range("AG13").Value = 18
$ws.Range("AH13").Value = 16
$ws.Range("AJ13").Value = 110
$ws.Range("AK13").Value = 48
$ws.Range("AL13").Value = 46
$ws.Range("AN13").Value = 38
$ws.Range("AO13").Value = 7.4
